$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document's body text:
# "Multiple forests can be synced to a single Azure AD directory. Azure AD connect per forest."
# followed by the _GoBack bookmark, then one trailing empty paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Multiple forests can be synced*") {
        $target = $p
    }
}

$r = $target.Range
$r.Collapse(0)

$lines = @(
    "",
    "What role must the account used for Azure AD Connect have? Global Admin Role",
    "What group membership must the on-premises account used for Azure AD Connect have? Enterprise Admin Group in the forest.",
    "",
    "",
    "",
    "Secure Resources",
    "",
    "",
    "Azure ACS (Azure Access Control Service)",
    "Depreciated but still supported",
    "Merged into Azure AD",
    "",
    "AZ FS (Federation Services)",
    "Syncs onpremise to Azure AD and use Azure AD to sync",
    "SSO with corporate network with policies that can be set.",
    "Use it when you need conditional access for onpremi/cloud and restricting sync of password with hash.",
    "Can setup now with Azure AD Connect with FS and web application proxy.",
    "web application proxy should be accessible to the internet.",
    "",
    "Enable `u{2013}PSRemoting `u{2013}force",
    "",
    "What role service, besides AD FS, is required when federating on-premises Active Directory with Azure AD when using Azure AD Connect with Windows Server 2012 R2? Web application porxy role and accessible to the internet.",
    ""
)

$blob = ($lines -join "`r") + "`r"
$r.Text = $blob

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
